# Updated cryptos list on Mon Aug  7 23:18:21 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values to the cryptos sheet,
# preserving each cell's original text representation (incl. trailing zero
# style numbers like '1.000' and the padded '  +0.12%  ' volume strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text value
$changes = [ordered]@{
    'D2' = '29.166.26'
    'E2' = '  +0.11%  '
    'D3' = '1.826.76'
    'E3' = '  -0.27%  '
    'D4' = '1.000'
    'E4' = '  -0.05%  '
    'D5' = '241.73'
    'E5' = '  -0.67%  '
    'D6' = '0.6206'
    'E6' = '  -0.82%  '
    'D7' = '1.001'
    'E7' = '  -0.26%  '
    'D8' = '0.07332'
    'E8' = '  -1.66%  '
    'D9' = '0.2894'
    'E9' = '  -1.07%  '
    'D10' = '22.97'
    'E10' = '  -1.39%  '
    'E11' = '  -0.30%  '
    'D12' = '1.823.96'
    'E12' = '  -0.64%  '
    'D13' = '4.960'
    'E13' = '  -1.01%  '
    'D14' = '0.6644'
    'E14' = '  -0.48%  '
    'D15' = '82.35'
    'E15' = '  -0.24%  '
    'D16' = '0.000008952'
    'E16' = '  -4.23%  '
    'D17' = '5.849'
    'E17' = '  -1.85%  '
    'D18' = '29.130.15'
    'E18' = '  -0.02%  '
    'D19' = '2.066.98'
    'E19' = '  -0.46%  '
    'D20' = '238.34'
    'E20' = '  +6.92%  '
    'D21' = '12.42'
    'E21' = '  -1.36%  '
    'D22' = '1.001'
    'E22' = '  -0.33%  '
    'D23' = '7.316'
    'E23' = '  +2.66%  '
    'D24' = '1.002'
    'E24' = '  -0.11%  '
    'D25' = '158.26'
    'E25' = '  -1.27%  '
    'D26' = '0.1422'
    'E26' = '  +2.16%  '
    'D27' = '8.472'
    'E27' = '  -0.26%  '
    'D28' = '17.66'
    'E28' = '  -1.36%  '
    'E29' = '  -0.58%  '
    'D30' = '0.05567'
    'E30' = '  -4.29%  '
    'D31' = '4.085'
    'E31' = '  -0.95%  '
    'D32' = '4.094'
    'E32' = '  -1.58%  '
    'D33' = '1.207'
    'E33' = '  -0.61%  '
    'D34' = '1.845'
    'E34' = '  +0.97%  '
    'D35' = '0.7347'
    'E35' = '  -0.55%  '
    'D36' = '1.131'
    'E36' = '  -0.67%  '
    'D37' = '2.624'
    'E37' = '  -1.94%  '
    'D38' = '2.845'
    'D39' = '1.211.23'
    'E39' = '  -1.36%  '
    'D40' = '0.01765'
    'E40' = '  -0.34%  '
    'D41' = '6.301'
    'E41' = '  -3.00%  '
    'D42' = '0.9159'
    'E42' = '  +2.44%  '
    'E43' = '  -0.23%  '
    'D44' = '101.66'
    'E44' = '  -0.48%  '
    'D45' = '1.973.12'
    'E45' = '  -0.26%  '
    'D46' = '64.71'
    'E46' = '  -1.83%  '
    'D47' = '0.5092'
    'E47' = '  +0.02%  '
    'E48' = '  -7.92%  '
    'D49' = '9.152'
    'E49' = '  +1.82%  '
    'D50' = '0.4024'
    'E50' = '  -0.89%  '
    'D51' = '0.05758'
    'E51' = '  -1.27%  '
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    # Force plain-text interpretation so numeric-looking strings (e.g. '1.000',
    # '29.166.26') are kept verbatim instead of being parsed into numbers, then
    # drop back to the unstyled 'Normal' style so no quote-prefix formatting
    # leaks into the saved cell style.
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
    $cell.Style = "Normal"
}
